$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 165.0625
$ws.Range("I9").Value = 171.8
$ws.Range("J9").Value = 64
$ws.Range("K9").Value = 171.8
$ws.Range("L9").Value = 64
$ws.Range("M9").Value = -2.800000000000011
$ws.Range("N9").Value = -402

$ws.Range("H19").Value = 1874.5385
$ws.Range("I19").Value = 1424.8334
$ws.Range("J19").Value = 2260
$ws.Range("K19").Value = 1424.8334
$ws.Range("L19").Value = 2260
$ws.Range("M19").Value = -1249.8334
$ws.Range("N19").Value = -2610

$ws.Range("H48").Value = 1438.8
$ws.Range("I48").Value = 548.5
$ws.Range("J48").Value = 5000
$ws.Range("K48").Value = 1645.5
$ws.Range("L48").Value = 15000
$ws.Range("M48").Value = -1353.5
$ws.Range("N48").Value = -15584

$ws.Range("H56").Value = 1438.8
$ws.Range("I56").Value = 548.5
$ws.Range("J56").Value = 5000
$ws.Range("K56").Value = 1645.5
$ws.Range("L56").Value = 15000
$ws.Range("M56").Value = -1111.5
$ws.Range("N56").Value = -16068

$ws.Range("H92").Value = 925.3333
$ws.Range("I92").Value = 925.3333
$ws.Range("K92").Value = 925.3333
$ws.Range("M92").Value = 322.6667

$ws.Range("H107").Value = 754
$ws.Range("I107").Value = 754
$ws.Range("K107").Value = 754
$ws.Range("M107").Value = 1166

$ws.Range("H113").Value = 83335690
$ws.Range("J113").Value = 3933.3333
$ws.Range("L113").Value = 3933.3333
$ws.Range("N113").Value = -10441.3333

$ws.Range("H127").Value = 666.6667
$ws.Range("I127").Value = 650
$ws.Range("J127").Value = 700
$ws.Range("K127").Value = 1950
$ws.Range("L127").Value = 2100
$ws.Range("M127").Value = 3010
$ws.Range("N127").Value = -12020

$ws.Range("H132").Value = 4545.222
$ws.Range("I132").Value = 4677.5386
$ws.Range("K132").Value = 14032.6158
$ws.Range("M132").Value = -11502.6158

$ws.Range("H135").Value = 1079.1538
$ws.Range("J135").Value = 93
$ws.Range("L135").Value = 837
$ws.Range("N135").Value = -5907

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 5754.116
$ws.Range("I32").Value = 5754.116
$ws.Range("J32").Value = 0
$ws.Range("K32").Value = 5754.116
$ws.Range("L32").Value = 0
$ws.Range("M32").Value = -5467.116
$ws.Range("N32").ClearContents()

$ws.Range("H74").Value = 5420.5415
$ws.Range("J74").Value = 7081.1816
$ws.Range("L74").Value = 7081.1816
$ws.Range("N74").Value = -8829.1816

$ws.Range("H77").Value = 5420.5415
$ws.Range("J77").Value = 7081.1816
$ws.Range("L77").Value = 35405.908
$ws.Range("N77").Value = -44141.908

$ws.Range("H88").Value = 2052.6667
$ws.Range("J88").Value = 2617.8
$ws.Range("L88").Value = 2617.8
$ws.Range("N88").Value = -3429.8

$ws.Range("H91").Value = 2052.6667
$ws.Range("J91").Value = 2617.8
$ws.Range("L91").Value = 2617.8
$ws.Range("N91").Value = -5425.8

$ws.Range("H97").Value = 387.7143
$ws.Range("I97").Value = 291.22223
$ws.Range("J97").Value = 966.6667
$ws.Range("K97").Value = 291.22223
$ws.Range("L97").Value = 966.6667
$ws.Range("M97").Value = 204.77777
$ws.Range("N97").Value = -1958.6667

$ws.Range("H101").Value = 46399.332
$ws.Range("J101").Value = 46399.332
$ws.Range("L101").Value = 46399.332
$ws.Range("N101").Value = -52889.332

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H74").Value = 0
$ws.Range("J74").Value = 0
$ws.Range("L74").Value = 0
$ws.Range("N74").ClearContents()

$ws.Range("H77").Value = 0
$ws.Range("J77").Value = 0
$ws.Range("L77").Value = 0
$ws.Range("N77").ClearContents()

$ws.Range("H80").Value = 439.57144
$ws.Range("J80").Value = 498.7
$ws.Range("L80").Value = 498.7
$ws.Range("N80").Value = -2494.7

$ws.Range("H83").Value = 439.57144
$ws.Range("J83").Value = 498.7
$ws.Range("L83").Value = 2493.5
$ws.Range("N83").Value = -12477.5

$ws.Range("H86").Value = 33336688
$ws.Range("I86").Value = 5500
$ws.Range("J86").Value = 50002280
$ws.Range("K86").Value = 5500
$ws.Range("L86").Value = 50002280
$ws.Range("M86").Value = -4377
$ws.Range("N86").Value = -50004526

$ws.Range("H89").Value = 33336688
$ws.Range("I89").Value = 5500
$ws.Range("J89").Value = 50002280
$ws.Range("K89").Value = 27500
$ws.Range("L89").Value = 250011400
$ws.Range("M89").Value = -21884
$ws.Range("N89").Value = -250022632

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 7216.9473
$ws.Range("J31").Value = 7724.727
$ws.Range("L31").Value = 7724.727
$ws.Range("N31").Value = -8314.726999999999

$ws.Range("H34").Value = 7216.9473
$ws.Range("J34").Value = 7724.727
$ws.Range("L34").Value = 7724.727
$ws.Range("N34").Value = -8128.727

$ws.Range("H94").Value = 5000
$ws.Range("I94").Value = 0
$ws.Range("K94").Value = 0
$ws.Range("M94").ClearContents()

$ws.Range("H99").Value = 5924.706
$ws.Range("J99").Value = 8283.166999999999
$ws.Range("L99").Value = 8283.166999999999
$ws.Range("N99").Value = -11279.167

$ws.Range("H107").Value = 1380.4482
$ws.Range("J107").Value = 2042.9231
$ws.Range("L107").Value = 2042.9231
$ws.Range("N107").Value = -5882.9231

$ws.Range("H111").Value = 74286.5
$ws.Range("J111").Value = 74286.5
$ws.Range("L111").Value = 74286.5
$ws.Range("N111").Value = -82466.5

$ws.Range("H112").Value = 64169
$ws.Range("J112").Value = 64169
$ws.Range("L112").Value = 64169
$ws.Range("N112").Value = -67123

$ws.Range("H126").Value = 5924.706
$ws.Range("J126").Value = 8283.166999999999
$ws.Range("L126").Value = 24849.501
$ws.Range("N126").Value = -29789.501

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H34").Value = 62822.188
$ws.Range("J34").Value = 143284.28
$ws.Range("L34").Value = 429852.84
$ws.Range("N34").Value = -430020.84

$ws.Range("H38").Value = 656.9286
$ws.Range("I38").Value = 26.666666
$ws.Range("J38").Value = 1129.625
$ws.Range("K38").Value = 79.99999800000001
$ws.Range("L38").Value = 3388.875
$ws.Range("M38").Value = 267.000002
$ws.Range("N38").Value = -4082.875

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H29").Value = 3017.75
$ws.Range("I29").Value = 3017.75
$ws.Range("J29").Value = 0
$ws.Range("K29").Value = 3017.75
$ws.Range("L29").Value = 0
$ws.Range("M29").Value = -2727.75
$ws.Range("N29").ClearContents()

$ws.Range("H43").Value = 5153.846
$ws.Range("I43").Value = 17000
$ws.Range("J43").Value = 3000
$ws.Range("K43").Value = 17000
$ws.Range("L43").Value = 3000
$ws.Range("M43").Value = -16849
$ws.Range("N43").Value = -3302

$ws.Range("H70").Value = 10914.833
$ws.Range("I70").Value = 3829.6667
$ws.Range("K70").Value = 3829.6667
$ws.Range("M70").Value = -3559.6667

$ws.Range("H73").Value = 10914.833
$ws.Range("I73").Value = 3829.6667
$ws.Range("K73").Value = 3829.6667
$ws.Range("M73").Value = -2893.6667

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2990.4
$ws.Range("I7").Value = 2933.7778
$ws.Range("K7").Value = 2933.7778
$ws.Range("M7").Value = -2821.7778

$ws.Range("H24").Value = 25000
$ws.Range("J24").Value = 25000
$ws.Range("L24").Value = 25000
$ws.Range("N24").Value = -25686

$ws.Range("H61").Value = 253252
$ws.Range("I61").Value = 337336
$ws.Range("J61").Value = 1000
$ws.Range("K61").Value = 337336
$ws.Range("L61").Value = 1000
$ws.Range("M61").Value = -337134
$ws.Range("N61").Value = -1404

$ws.Range("H113").Value = 253252
$ws.Range("I113").Value = 337336
$ws.Range("J113").Value = 1000
$ws.Range("K113").Value = 337336
$ws.Range("L113").Value = 1000
$ws.Range("M113").Value = -335166
$ws.Range("N113").Value = -5340

$ws.Range("H126").Value = 2990.4
$ws.Range("I126").Value = 2933.7778
$ws.Range("K126").Value = 8801.3334
$ws.Range("M126").Value = -6331.3334

$ws.Range("H132").Value = 8597.77
$ws.Range("I132").Value = 8787.951999999999
$ws.Range("K132").Value = 26363.856
$ws.Range("M132").Value = -23833.856

$ws.Range("H140").Value = 94158
$ws.Range("J140").Value = 94158
$ws.Range("L140").Value = 94158
$ws.Range("N140").Value = -104518

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H26").Value = 20000
$ws.Range("I26").Value = 20000
$ws.Range("K26").Value = 20000
$ws.Range("M26").Value = -19707

$ws.Range("H27").Value = 0
$ws.Range("J27").Value = 0
$ws.Range("L27").Value = 0
$ws.Range("N27").ClearContents()

$ws.Range("H31").Value = 0
$ws.Range("J31").Value = 0
$ws.Range("L31").Value = 0
$ws.Range("N31").ClearContents()

$ws.Range("H125").Value = 64490
$ws.Range("J125").Value = 64490
$ws.Range("L125").Value = 64490
$ws.Range("N125").Value = -74330

$ws.Range("H136").Value = 8096.913
$ws.Range("I136").Value = 6951.8125
$ws.Range("K136").Value = 20855.4375
$ws.Range("M136").Value = -18305.4375
